# "Add price for all stocks"
# Adds a new "Price" column (N) with the header in N1 and per-row stock
# prices in N7:N50 (N2:N6 are left blank, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("N1").Value = "Price"

# Give the whole N2:N50 block the same "black Calibri" font used by the
# bulk of the pasted-in price data (this also covers N7:N50 below; the
# currency number format + the one-off font override for N7 are applied
# on top of it afterwards).
$ws.Range("N2:N50").Font.Color = 0

# Price values (row 7 is the first data row carried over from the paste).
$ws.Range("N7").Value = 27.11
$ws.Range("N8").Value = 29.26
$ws.Range("N9").Value = 30.55
$ws.Range("N10").Value = 39.72
$ws.Range("N11").Value = 41.87
$ws.Range("N12").Value = 32.26
$ws.Range("N13").Value = 42.12
$ws.Range("N14").Value = 47.29
$ws.Range("N15").Value = 35.93
$ws.Range("N16").Value = 40.61
$ws.Range("N17").Value = 34.14
$ws.Range("N18").Value = 29.08
$ws.Range("N19").Value = 18.28
$ws.Range("N20").Value = 20.46
$ws.Range("N21").Value = 17.23
$ws.Range("N22").Value = 8.99
$ws.Range("N23").Value = 15.72
$ws.Range("N24").Value = 12.06
$ws.Range("N25").Value = 27.97
$ws.Range("N26").Value = 23.96
$ws.Range("N27").Value = 23.25
$ws.Range("N28").Value = 30.61
$ws.Range("N29").Value = 30.35
$ws.Range("N30").Value = 37.81
$ws.Range("N31").Value = 44.94
$ws.Range("N32").Value = 39.47
$ws.Range("N33").Value = 39.55
$ws.Range("N34").Value = 42.23
$ws.Range("N35").Value = 49.48
$ws.Range("N36").Value = 48.79
$ws.Range("N37").Value = 56.93
$ws.Range("N38").Value = 59.48
$ws.Range("N39").Value = 49.87
$ws.Range("N40").Value = 54.24
$ws.Range("N41").Value = 35.33
$ws.Range("N42").Value = 47.94
$ws.Range("N43").Value = 46.47
$ws.Range("N44").Value = 66.79
$ws.Range("N45").Value = 67.31
$ws.Range("N46").Value = 75.24
$ws.Range("N47").Value = 101.32
$ws.Range("N48").Value = 87.83
$ws.Range("N49").Value = 78.5
$ws.Range("N50").Value = 90.85

# Currency number format for all the price cells.
$ws.Range("N7:N50").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# The very first pasted price (N7) kept the source document's font
# (Aptos Narrow) instead of picking up the workbook's Calibri.
$ws.Range("N7").Font.Name = "Aptos Narrow"

# Selection left on the pasted column, like Excel does right after a paste.
$ws.Range("N2:N50").Select()
